# chore: update Sheets via scheduled runner
# Refreshes cached profit-margin figures (columns H-N) on several rows
# across the ALC/ARM/BSM/CRP/CUL/GSM/LTW/WVR sheets.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H19").Value = 447.72726
$ws.Range("J19").Value = 420.66666
$ws.Range("L19").Value = 420.66666
$ws.Range("N19").Value = -770.66666

$ws.Range("H49").Value = 300
$ws.Range("I49").Value = 200
$ws.Range("K49").Value = 600
$ws.Range("M49").Value = -464

$ws.Range("H74").Value = 3455.1904
$ws.Range("I74").Value = 3349.9375
$ws.Range("J74").Value = 3792
$ws.Range("K74").Value = 3349.9375
$ws.Range("L74").Value = 3792
$ws.Range("M74").Value = -2413.9375
$ws.Range("N74").Value = -5664

$ws.Range("H77").Value = 3455.1904
$ws.Range("I77").Value = 3349.9375
$ws.Range("J77").Value = 3792
$ws.Range("K77").Value = 16749.6875
$ws.Range("L77").Value = 18960
$ws.Range("M77").Value = -12069.6875
$ws.Range("N77").Value = -28320

$ws.Range("H107").Value = 389.63333
$ws.Range("I107").Value = 408.94736
$ws.Range("J107").Value = 356.27274
$ws.Range("K107").Value = 408.94736
$ws.Range("L107").Value = 356.27274
$ws.Range("M107").Value = 1511.05264
$ws.Range("N107").Value = -4196.27274

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1461.5
$ws.Range("I2").Value = 1423.7273
$ws.Range("J2").Value = 1600
$ws.Range("K2").Value = 1423.7273
$ws.Range("L2").Value = 1600
$ws.Range("M2").Value = -1310.7273
$ws.Range("N2").Value = -1826

$ws.Range("H32").Value = 4184.84
$ws.Range("I32").Value = 3122.8333
$ws.Range("K32").Value = 3122.8333
$ws.Range("M32").Value = -2835.8333

$ws.Range("H61").Value = 942.75
$ws.Range("I61").Value = 942.75
$ws.Range("K61").Value = 942.75
$ws.Range("M61").Value = -730.75

$ws.Range("H116").Value = 1461.5
$ws.Range("I116").Value = 1423.7273
$ws.Range("J116").Value = 1600
$ws.Range("K116").Value = 1423.7273
$ws.Range("L116").Value = 1600
$ws.Range("M116").Value = 870.2727
$ws.Range("N116").Value = -6188

$ws.Range("H132").Value = 1648
$ws.Range("I132").Value = 497.2
$ws.Range("J132").Value = 3566
$ws.Range("K132").Value = 1491.6
$ws.Range("L132").Value = 10698
$ws.Range("M132").Value = 1038.4
$ws.Range("N132").Value = -15758

$ws.Range("H136").Value = 942.75
$ws.Range("I136").Value = 942.75
$ws.Range("K136").Value = 2828.25
$ws.Range("M136").Value = -278.25

$ws.Range("H139").Value = 124000
$ws.Range("I139").Value = 0
$ws.Range("J139").Value = 124000
$ws.Range("K139").Value = 0
$ws.Range("L139").Value = 124000
$ws.Range("M139").ClearContents()
$ws.Range("N139").Value = -134280

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1461.5
$ws.Range("I3").Value = 1423.7273
$ws.Range("J3").Value = 1600
$ws.Range("K3").Value = 1423.7273
$ws.Range("L3").Value = 1600
$ws.Range("M3").Value = -1309.7273
$ws.Range("N3").Value = -1828

$ws.Range("H105").Value = 4243.5
$ws.Range("J105").Value = 955.5
$ws.Range("L105").Value = 955.5
$ws.Range("N105").Value = -4449.5

$ws.Range("H107").Value = 691.1111
$ws.Range("I107").Value = 652.5
$ws.Range("J107").Value = 1000
$ws.Range("K107").Value = 652.5
$ws.Range("L107").Value = 1000
$ws.Range("M107").Value = 1267.5
$ws.Range("N107").Value = -4840

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 30383.223
$ws.Range("I31").Value = 44754.695
$ws.Range("J31").Value = 4956.769
$ws.Range("K31").Value = 44754.695
$ws.Range("L31").Value = 4956.769
$ws.Range("M31").Value = -44459.695
$ws.Range("N31").Value = -5546.769

$ws.Range("H34").Value = 30383.223
$ws.Range("I34").Value = 44754.695
$ws.Range("J34").Value = 4956.769
$ws.Range("K34").Value = 44754.695
$ws.Range("L34").Value = 4956.769
$ws.Range("M34").Value = -44552.695
$ws.Range("N34").Value = -5360.769

$ws.Range("H58").Value = 1247.3914
$ws.Range("I58").Value = 1268.4722
$ws.Range("J58").Value = 1171.5
$ws.Range("K58").Value = 1268.4722
$ws.Range("L58").Value = 1171.5
$ws.Range("M58").Value = -1065.4722
$ws.Range("N58").Value = -1577.5

$ws.Range("H62").Value = 2720.6
$ws.Range("I62").Value = 2614.2856
$ws.Range("J62").Value = 2968.6667
$ws.Range("K62").Value = 2614.2856
$ws.Range("L62").Value = 2968.6667
$ws.Range("M62").Value = -1990.2856
$ws.Range("N62").Value = -4216.6667

$ws.Range("H65").Value = 2720.6
$ws.Range("I65").Value = 2614.2856
$ws.Range("J65").Value = 2968.6667
$ws.Range("K65").Value = 13071.428
$ws.Range("L65").Value = 14843.3335
$ws.Range("M65").Value = -9951.428
$ws.Range("N65").Value = -21083.3335

$ws.Range("H107").Value = 1404.4667
$ws.Range("I107").Value = 1630.3334
$ws.Range("J107").Value = 501
$ws.Range("K107").Value = 1630.3334
$ws.Range("L107").Value = 501
$ws.Range("M107").Value = 289.6666
$ws.Range("N107").Value = -4341

$ws.Range("H122").Value = 3423.2307
$ws.Range("I122").Value = 3625.1667
$ws.Range("K122").Value = 10875.5001
$ws.Range("M122").Value = -8425.500100000001

$ws.Range("H136").Value = 1247.3914
$ws.Range("I136").Value = 1268.4722
$ws.Range("J136").Value = 1171.5
$ws.Range("K136").Value = 3805.4166
$ws.Range("L136").Value = 3514.5
$ws.Range("M136").Value = -1255.4166
$ws.Range("N136").Value = -8614.5

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H131").Value = 1327048
$ws.Range("J131").Value = 1603438.4
$ws.Range("L131").Value = 4810315.199999999
$ws.Range("N131").Value = -4820395.199999999

$ws.Range("H138").Value = 682.9167
$ws.Range("I138").Value = 635.9091
$ws.Range("J138").Value = 1200
$ws.Range("K138").Value = 1907.7273
$ws.Range("L138").Value = 3600
$ws.Range("M138").Value = 3232.2727
$ws.Range("N138").Value = -13880

$ws.Range("H140").Value = 2334.95
$ws.Range("I140").Value = 2336.9
$ws.Range("J140").Value = 2333
$ws.Range("K140").Value = 7010.700000000001
$ws.Range("L140").Value = 6999
$ws.Range("M140").Value = -1830.700000000001
$ws.Range("N140").Value = -17359

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 210
$ws.Range("I2").Value = 170
$ws.Range("J2").Value = 270
$ws.Range("K2").Value = 170
$ws.Range("L2").Value = 270
$ws.Range("M2").Value = -57
$ws.Range("N2").Value = -496

$ws.Range("H122").Value = 2901.75
$ws.Range("I122").Value = 3007
$ws.Range("J122").Value = 2866.6667
$ws.Range("K122").Value = 9021
$ws.Range("L122").Value = 8600.000100000001
$ws.Range("M122").Value = -6571
$ws.Range("N122").Value = -13500.0001

$ws.Range("H132").Value = 4001.6667
$ws.Range("I132").Value = 3004
$ws.Range("J132").Value = 4999.3335
$ws.Range("K132").Value = 9012
$ws.Range("L132").Value = 14998.0005
$ws.Range("M132").Value = -6482
$ws.Range("N132").Value = -20058.0005

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H122").Value = 6411.273
$ws.Range("I122").Value = 11126
$ws.Range("J122").Value = 3717.1428
$ws.Range("K122").Value = 33378
$ws.Range("L122").Value = 11151.4284
$ws.Range("M122").Value = -30928
$ws.Range("N122").Value = -16051.4284

$ws.Range("H136").Value = 2357.5508
$ws.Range("I136").Value = 1715.381
$ws.Range("J136").Value = 3356.4814
$ws.Range("K136").Value = 5146.143
$ws.Range("L136").Value = 10069.4442
$ws.Range("M136").Value = -2596.143
$ws.Range("N136").Value = -15169.4442

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H118").Value = 0
$ws.Range("J118").Value = 0
$ws.Range("L118").Value = 0
$ws.Range("N118").ClearContents()

$ws.Range("H122").Value = 2002236.8
$ws.Range("I122").Value = 2501796
$ws.Range("K122").Value = 7505388
$ws.Range("M122").Value = -7502938
